$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7, 8 and 9 - those invoices are no longer part of the export
$ws.Rows.Item(7).Resize(3).Delete()

# Dates are stored as plain text (dd/mm/yyyy), keep them as text rather than
# letting Excel auto-convert them into date serials. Force the "@" text
# format while entering the value, then clear the formatting again so the
# cell's style matches the rest of the untouched data cells.
$ws.Range("C2:C6").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "TW1_033791184"
$ws.Range("C2").Value = "04/10/2024"
$ws.Range("C2").ClearFormats()
$ws.Range("D2").Value = "17, Rue de Téhéran, 75008 PARIS- 8EME"
$ws.Range("E2").Value = "28, Rue Petit, 92110 CLICHY"
$ws.Range("F2").Value = "20,00 €"
$ws.Range("G2").Value = "TW1_033791184_004089382_00304015690.pdf"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "TW1_033844533"
$ws.Range("C3").Value = "07/10/2024"
$ws.Range("C3").ClearFormats()
$ws.Range("D3").Value = "17, Rue de Téhéran, 75008 PARIS- 8EME"
$ws.Range("E3").Value = "28, Rue Petit, 92110 CLICHY"
$ws.Range("F3").Value = "29,50 €"
$ws.Range("G3").Value = "TW1_033844533_004089382_00304130605.pdf"

# Row 4
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "TW1_033856935"
$ws.Range("C4").Value = "07/10/2024"
$ws.Range("C4").ClearFormats()
$ws.Range("D4").Value = "28, Rue Petit, 92110 CLICHY"
$ws.Range("E4").Value = "17, Rue de Téhéran, 75008 PARIS- 8EME"
$ws.Range("F4").Value = "23,10 €"
$ws.Range("G4").Value = "TW1_033856935_004089382_00304167026.pdf"

# Row 5
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "TW1_033894205"
$ws.Range("C5").Value = "09/10/2024"
$ws.Range("C5").ClearFormats()
$ws.Range("D5").Value = "17, Rue de Téhéran, 75008 PARIS- 8EME"
$ws.Range("E5").Value = "28, Rue Petit, 92110 CLICHY"
$ws.Range("F5").Value = "22,40 €"
$ws.Range("G5").Value = "TW1_033894205_004089382_00304261950.pdf"

# Row 6
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "TW1_033901495"
$ws.Range("C6").Value = "09/10/2024"
$ws.Range("C6").ClearFormats()
$ws.Range("D6").Value = "28, Rue Petit, 92110 CLICHY"
$ws.Range("E6").Value = "15, Rue Charles Duflos, 92270 BOIS-COLOMBES"
$ws.Range("F6").Value = "21,10 €"
$ws.Range("G6").Value = "TW1_033901495_004089382_00304282644.pdf"
